$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a Date column (H) - header is a space " " string for each row 2-7,
# and H1 holds a date value (2015-01-01, serial 42005)
$ws.Range("H1").NumberFormat = "mm-dd-yy"
$ws.Range("H1").Value = Get-Date -Year 2015 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0

$ws.Range("H2").Value = " "
$ws.Range("H3").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("H5").Value = " "
$ws.Range("H6").Value = " "
$ws.Range("H7").Value = " "

# Update selection to match the new active cell
$ws.Range("H7").Select()
